$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data rows 2-14, columns A (Meta), B (Venda), C (Pecas)
$data = @(
    @(5000, 5000, 30),
    @(5000, 6000, 10),
    @(5000, 6000, 10),
    @(5000, 6000, 10),
    @(5000, 6000, 10),
    @(5000, 6000, 10),
    @(5000, 5000, 700),
    @(500,  5000, 50),
    @(5000, 6000, 10),
    @(5000, 6000, 10),
    @(5000, 6000, 10),
    @(5000, 6000, 10),
    @(5000, 6000, 10)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
